# "Generate Report for Handoff"
#
# The localization-status report just finished a handoff run: the
# zh-cn / de-de rows flip from "In Translation" to "Ready for handoff"
# and the handoff timestamps on the Overview + per-language sheets move
# forward a couple of minutes. Excel also re-autosized the (now wider)
# status/date columns on each sheet to fit the new text.

$wb = $excel.ActiveWorkbook

# Column E/F on Overview (and the Status column on the language sheets)
# went from ~13.41 chars to ~17.22 chars wide to fit "Ready for handoff".
# This engine's ColumnWidth setter snaps to 1/6-character increments, so
# 16.3333 is the input that lands on the closest reachable stored width.
$autoFitStatusColumnWidth = 16.333333333333336

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-31 08:23:26"
$wsOverview.Columns.Item(5).ColumnWidth = $autoFitStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $autoFitStatusColumnWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-31 08:23:15"
$wsZhCn.Columns.Item(3).ColumnWidth = $autoFitStatusColumnWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-31 08:23:26"
$wsDeDe.Columns.Item(3).ColumnWidth = $autoFitStatusColumnWidth
